# Add a new worksheet "FTNC_Average_Demand101" after the existing
# "FTNC_Average_Demand10" sheet, as a copy of it (so headers, styles,
# sheetPr and page margins all carry over), then overwrite the numeric
# results row with the new figures.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("FTNC_Average_Demand10")

# Copy the source sheet and place the copy immediately after it.
$src.Copy($null, $src)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FTNC_Average_Demand101"

# Update the figures (modify the figure display) on the new sheet.
$newSheet.Range("B2").Value = 2455.274984722344
$newSheet.Range("C2").Value = 13095.80869023837
$newSheet.Range("D2").Value = 630.3256680479443
$newSheet.Range("E2").Value = 31.57236976815668
$newSheet.Range("F2").Value = 16212.98171277687
